$wb = $excel.ActiveWorkbook

# --- Table_Names sheet ---
$wsTables = $wb.Worksheets.Item("Table_Names")
$wsTables.Range("A1").Value = "TestTable1"
$wsTables.Range("A2").Value = "TestTable2"
$wsTables.Range("A3").Value = "TestTable3"

# --- Field_Names sheet ---
$wsFields = $wb.Worksheets.Item("Field_Names")
$wsFields.Range("A1").Value = "TestValue1.1"
$wsFields.Range("A2").Value = "TestValue1.2"
$wsFields.Range("A3").Value = "TestValue1.3"
$wsFields.Range("A4").Value = "TestValue2.1"
$wsFields.Range("A5").Value = "TestValue2.2"
$wsFields.Range("A6").Value = "TestValue3.1"
$wsFields.Range("A7").Value = "TestValue3.2"
$wsFields.Range("A8").Value = "TestValue3.3"
$wsFields.Range("A9").Value = "TestValue3.4"
